$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 171, shifting existing rows 171:199 down to 172:200
$ws.Rows.Item(171).Insert()

# Populate the new row 171 with the new record
$ws.Range("A171").Value = 5
$ws.Range("B171").Value = "Macroferia Regional de Talca"
$ws.Range("C171").Value = "Maule"
$ws.Range("D171").Value = 44505
$ws.Range("D171").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E171").Value = 7
$ws.Range("F171").Value = 100114014
$ws.Range("G171").Value = "Betarraga"
$ws.Range("H171").Value = "Sin especificar"
$ws.Range("I171").Value = "Primera"
$ws.Range("J171").Value = 5000
$ws.Range("K171").Value = 650
$ws.Range("L171").Value = 650
$ws.Range("M171").Value = 650
$ws.Range("N171").Value = "$/paquete 5 unidades"
$ws.Range("O171").Value = "Región del Maule"
$ws.Range("P171").Value = 130
$ws.Range("Q171").Value = 5
$ws.Range("R171").Value = "Hortaliza"
